# Insert 4 new data rows right above the existing row 330, pushing the
# current rows 330..409 down to 334..413 (dimension grows from T409 to T413).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("330:333").Insert()

# Shared/constant column values for this dataset block.
$marketId   = 7
$market     = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$catId      = 100103004
$categoria  = "Durazno"
$origen     = "Región de O'Higgins"

# New rows to populate (A..T), matching the after-state of the diff.
$newRows = @(
    @{ Row = 330; Fecha = 44985; Variedad = "Carson";   Calidad = "Especial"; Vol = 50; PMin = 15000; PMax = 15000; PProm = 15000; Unidad = "`$/caja 16 kilos empedrada"; PKg = 938; KgUnidad = 16 },
    @{ Row = 331; Fecha = 44985; Variedad = "Carson";   Calidad = "Primera";  Vol = 50; PMin = 13000; PMax = 13000; PProm = 13000; Unidad = "`$/caja 16 kilos empedrada"; PKg = 812; KgUnidad = 16 },
    @{ Row = 332; Fecha = 44985; Variedad = "Carson";   Calidad = "Segunda";  Vol = 50; PMin = 11000; PMax = 11000; PProm = 11000; Unidad = "`$/caja 16 kilos empedrada"; PKg = 688; KgUnidad = 16 },
    @{ Row = 333; Fecha = 44985; Variedad = "Kakamas";  Calidad = "Primera";  Vol = 60; PMin = 14000; PMax = 14000; PProm = 14000; Unidad = "`$/caja 16 kilos empedrada"; PKg = 875; KgUnidad = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $marketId
    $ws.Cells.Item($row, 2).Value2  = $market
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $catId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $r.Variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Vol
    $ws.Cells.Item($row, 14).Value2 = $r.PMin
    $ws.Cells.Item($row, 15).Value2 = $r.PMax
    $ws.Cells.Item($row, 16).Value2 = $r.PProm
    $ws.Cells.Item($row, 17).Value2 = $r.Unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PKg
    $ws.Cells.Item($row, 20).Value2 = $r.KgUnidad
}
